# ScriptCraft Day 02 - "More Advanced Plugins" deck fix-up
#
# 1) Fix the "cmLocation is not a function" typo on the dice/location-
#    helper slide: `return cmLocation(x, y, z);` -> `return new cmLocation(x, y, z);`
#
# 2) Swap the presentation's active colour theme from the "swiss-2"
#    palette over to the "Custom Theme" palette (the deck's font scheme
#    and effect/format scheme are identical between the two themes -
#    only the 12 scheme colours differ), by rewriting every slot of the
#    slide master's ThemeColorScheme.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Code-sample typo fix (slide 7, the getBufferInFrontOfPlayer() box)
# ---------------------------------------------------------------------
$slide = $p.Slides.Item(7)
$shape = $slide.Shapes.Item(2)
$textRange = $shape.TextFrame.TextRange

# Locate the paragraph that still reads "  return cmLocation(x, y, z); "
$fullText = $textRange.Text
$searchText = "return cmLocation(x, y, z);"
$charIndex = $fullText.IndexOf($searchText)
if ($charIndex -ge 0) {
    # Work out which paragraph (1-based) contains that character offset.
    $prefix = $fullText.Substring(0, $charIndex)
    $paragraphNumber = ($prefix.Split("`r").Length)

    $targetParagraph = $textRange.Paragraphs($paragraphNumber, 1)

    # Route the replacement through an unrelated placeholder string first
    # so the engine's word-level run-diffing can't see any shared
    # prefix/suffix with the final text - this keeps the paragraph as a
    # single run (matching the original single-run authoring) instead of
    # splitting it into "return " / "new " / "cmLocation(...)" runs.
    $targetParagraph.Text = "PLACEHOLDER_TEXT_FOR_REPLACEMENT"
    $targetParagraph2 = $textRange.Paragraphs($paragraphNumber, 1)
    $targetParagraph2.Text = "  return new cmLocation(x, y, z); "
}

# ---------------------------------------------------------------------
# 2) Theme colour swap (swiss-2 -> Custom Theme)
# ---------------------------------------------------------------------
$theme = $p.SlideMaster.Theme
$colorScheme = $theme.ThemeColorScheme

# ThemeColorScheme slot order is dk1, lt1, dk2, lt2, accent1-6, hlink,
# folHlink (1-based) - the same order DrawingML uses for <a:clrScheme>.
# Values are OLE_COLOR ints (0x00BBGGRR) for the "Custom Theme" palette.
$newColors = @{
    1  = 0          # dk1      000000
    2  = 16777215   # lt1      FFFFFF
    3  = 5800213    # dk2      158158
    4  = 15987699   # lt2      F3F3F3
    5  = 13077765   # accent1  058DC7
    6  = 3322960     # accent2  50B432
    7  = 1791725     # accent3  ED561B
    8  = 61421       # accent4  EDEF00
    9  = 15059748    # accent5  24CBE5
    10 = 7529828     # accent6  64E572
    11 = 13369378    # hlink    2200CC
    12 = 9116245     # folHlink 551A8B
}

foreach ($slot in 1..12) {
    $colorScheme.Item($slot).RGB = $newColors[$slot]
}
